$d = $word.ActiveDocument

$insertText = "SaaS implementation and management, "
$marker = "creating SAML sign on integrations, "

# ---------------------------------------------------------------------
# Helper: insert $Text right at $Pos, then toggle Bold off/on across
# exactly the inserted span so Word is forced to give it its own <w:r>
# instead of silently merging it into a neighbour run. $preserveFrom,
# when >= 0, is an offset further back marking an existing run boundary
# that must survive too - InsertAfter can re-coalesce same-formatted
# runs across the whole paragraph, so any boundary we care about must be
# re-asserted *after* the insertion (doing it beforehand does not stick).
# ---------------------------------------------------------------------
function Insert-Split($Doc, $Pos, $Text, $PreserveFrom) {
    $insertPoint = $Doc.Range($Pos, $Pos)
    $insertPoint.InsertAfter($Text)

    if ($PreserveFrom -ge 0) {
        $preserveRng = $Doc.Range($PreserveFrom, $Pos)
        $preserveRng.Bold = 1
        $preserveRng.Bold = 0
    }

    $newRng = $Doc.Range($Pos, $Pos + $Text.Length)
    $newRng.Bold = 1
    $newRng.Bold = 0
}

# --- Locate occurrence 1: a single run covering the whole sentence ---
$fullSentence1 = "Responsible for onboarding new vendors and affiliates into the environment, including setting up site-to-site VPN tunnels, " + $marker + `
    "user and group management, and peripheral device configuration assistance."

$rng1 = $d.Content
$found1 = $rng1.Find.Execute($fullSentence1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find occurrence 1 (single-run) sentence"
}
$occ1End = $rng1.End
$markerRng1 = $d.Range($rng1.Start, $occ1End)
[void]$markerRng1.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos1 = $markerRng1.End

# --- Locate occurrence 2: split across two runs in the original file -
#   run A: "Responsible for onboarding new vendors and affiliates into the environment"
#   run B: ", including setting up site-to-site VPN tunnels, creating SAML
#           sign on integrations, user and group management, and
#           peripheral device configuration assistance."
# Search *after* occurrence 1 so we don't match its (textually identical)
# tail again.
$runBText = ", including setting up site-to-site VPN tunnels, " + $marker + `
    "user and group management, and peripheral device configuration assistance."

$rng2 = $d.Range($occ1End, $d.Content.End)
$found2 = $rng2.Find.Execute($runBText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find occurrence 2 (split-run) sentence"
}
$runBStart = $rng2.Start
$markerRng2 = $d.Range($runBStart, $rng2.End)
[void]$markerRng2.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos2 = $markerRng2.End

# Apply edits from the later occurrence back to the earlier one so that
# earlier offsets are never invalidated by a later insertion.
Insert-Split $d $pos2 $insertText $runBStart
Insert-Split $d $pos1 $insertText -1

Write-Host "Done"
